$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "eta" column (D) is being split into two columns: "eta_min" and "eta_max".
# Insert a new blank column at E so everything from the old D onward shifts
# right by one (old D stays put for now, old E..L become F..M).
$ws.Columns("E:E").Insert()

# Re-purpose column D ("eta" -> "eta_min") and fill in the new column E ("eta_max").
$ws.Range("D1").Value = "eta_min"
$ws.Range("E1").Value = "eta_max"

$ws.Range("D2").Value = -0.17
$ws.Range("E2").Value = 0.17

$ws.Range("D3").Value = -0.17
$ws.Range("E3").Value = 0.17

# Match the author's final selection state.
$ws.Range("G11").Select()
